$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item("Textfeld 5")
$tr = $shape.TextFrame.TextRange
$para3 = $tr.Paragraphs(3, 1)
$para3.Text = " "
$para3b = $tr.Paragraphs(3, 1)
$para3b.Text = "Bis zum 3. Dezember  @TelekomCareer &    @wir_sind_die_onsite Instagram!"
